$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '94.249.22'
$ws.Range("E2").Value = '  +2.25%  '
$ws.Range("D3").Value = '3.072.46'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.97'
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '610.42'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.09'
$ws.Range("E7").Value = '  +0.76%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.377'
$ws.Range("E8").Value = '  -3.28%  '
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.801'
$ws.Range("E10").Value = '  +9.23%  '
$ws.Range("D11").Value = '3.071.45'
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("E12").Value = '  -1.77%  '
$ws.Range("D13").Value = '94.011.15'
$ws.Range("E13").Value = '  +1.96%  '
$ws.Range("E14").Value = '  -2.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '33.70'
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.33'
$ws.Range("E16").Value = '  -1.62%  '
$ws.Range("D17").Value = '3.643.19'
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("D18").Value = '3.059.05'
$ws.Range("E18").Value = '  -1.51%  '
$ws.Range("E19").Value = '  -5.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.35'
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("E21").Value = '  -0.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '438.38'
$ws.Range("E22").Value = '  -1.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.82'
$ws.Range("E23").Value = '  -5.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000189'
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("E25").Value = '  +6.46%  '
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '84.71'
$ws.Range("E27").Value = '  -1.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.84'
$ws.Range("E28").Value = '  +2.16%  '
$ws.Range("D29").Value = '3.236.29'
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.246'
$ws.Range("E31").Value = '  +9.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.178'
$ws.Range("E32").Value = '  +6.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.123'
$ws.Range("E33").Value = '  -6.74%  '
$ws.Range("E34").Value = '  -0.37%  '
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.979'
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.81'
$ws.Range("E36").Value = '  -0.60%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.154'
$ws.Range("E37").Value = '  -2.54%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.37'
$ws.Range("E38").Value = '  -1.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.89'
$ws.Range("E39").Value = '  +0.21%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '483.41'
$ws.Range("E40").Value = '  +0.78%  '
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '24.03'
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("B42").Value = 'MantraDAO'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.74'
$ws.Range("E42").Value = '  -3.71%  '
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.28'
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.436'
$ws.Range("E44").Value = '  +1.70%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.11'
$ws.Range("E46").Value = '  -5.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '161.47'
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("E48").Value = '  -1.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.81'
$ws.Range("E49").Value = '  -2.48%  '
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("E51").Value = '  +0.18%  '
